$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the split runs ("Name" / ", " / "Number") in each NAICS list
#    item back into a single run per paragraph. Doing a Find/Replace over
#    the full text of each item causes Word to rewrite the paragraph with
#    one merged run, which is exactly what the target markup needs.
# ---------------------------------------------------------------------------
$naicsLines = @(
    "Beauty Salons, 110895",
    "Beef Cattle Ranching and Farming, 50307",
    "Barber Shops, 48291",
    "All Other Personal Services, 45761",
    "General Freight Trucking, Local, 38820",
    "Offices of Real Estate Agents and Brokers, 38485",
    "Corn Farming, 37981",
    "General Freight Trucking, Long-Distance, Truckload, 36803",
    "Taxi Service, 36172",
    "Residential Remodelers, 29735"
)

foreach ($line in $naicsLines) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $found = $find.Execute($line, $false, $false, $false, $false, $false, $true, 1, $false, $line, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find/merge run for: $line"
    }
}

# ---------------------------------------------------------------------------
# 2) Append two brand-new plain (non-list) paragraphs at the very end of the
#    document, after the last "Compare the represented..." bullet and
#    before the closing sectPr. InsertXML is used (instead of
#    InsertParagraphAfter + Style=) so the new paragraphs come out as plain
#    <w:p> elements with no inherited ListParagraph/numPr formatting.
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$newParasXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Lendistry</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> is a character, but stands for something larger &#8211; banks with disproportionally rates of approval but not disbursing &#8211; why is that&gt;? what are the circumstances? who is affected? what happened to their businesses?</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t>break it down by industry, by location, look for specific individuals to talk to and do reporting</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$endRange.InsertXML($newParasXml)

Write-Host "edit complete"
